$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay as text so numeric-looking strings
# like "214.58" or "0.0626" are not silently coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.012.83'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '1.619.52'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '214.58'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.88%  '
$ws.Range('D9').Value = '0.0626'
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').Value = '1.618.22'
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = '64.66'
$ws.Range('E15').Value = '  -3.66%  '
$ws.Range('D16').Value = '27.003.07'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '0.0₃0747'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').Value = '214.28'
$ws.Range('E18').Value = '  -2.64%  '
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('D21').Value = '4.36'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('E22').Value = '  -5.77%  '
$ws.Range('D23').Value = '9.04'
$ws.Range('E23').Value = '  -1.39%  '
$ws.Range('D24').Value = '148.28'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '7.42'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  -1.86%  '
$ws.Range('D28').Value = '15.57'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('D29').Value = '0.0515'
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('E30').Value = '  -1.07%  '
$ws.Range('D31').Value = '3.36'
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('D32').Value = '0.748'
$ws.Range('E32').Value = '  +35.27%  '
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').Value = '1.347.39'
$ws.Range('E34').Value = '  +3.08%  '
$ws.Range('E35').Value = '  -1.03%  '
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('D38').Value = '0.849'
$ws.Range('E38').Value = '  -1.32%  '
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').Value = '0.801'
$ws.Range('E40').Value = '  -1.13%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').Value = '65.28'
$ws.Range('E42').Value = '  +5.05%  '
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = '1.757.07'
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('D45').Value = '89.72'
$ws.Range('E45').Value = '  -2.66%  '
$ws.Range('D46').Value = '0.867'
$ws.Range('E46').Value = '  +29.37%  '
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('D50').Value = '0.0999'
$ws.Range('E50').Value = '  +3.42%  '
$ws.Range('D51').Value = '7.69'
$ws.Range('E51').Value = '  +0.08%  '

# Restore default (General) formatting on column D so no residual
# text-format style is left on the cells themselves.
$ws.Range("D2:D51").ClearFormats()
